{"js": "const body = context.document.body;\n\n// Every \"old\" equation string below occurs exactly once in the document,\n// so a scoped search + single-hit replace safely retargets each cell in turn.\nconst pairs = [\n  [\"13\u00d721=273\", \"42\u00d729=1218\"],\n  [\"72\u00d770=5040\", \"41\u00d714=574\"],\n  [\"86\u00d735=3010\", \"20\u00d790=1800\"],\n  [\"24\u00d739=936\", \"11\u00d714=154\"],\n  [\"24\u00d795=2280\", \"47\u00d762=2914\"],\n  [\"88\u00d757=5016\", \"78\u00d744=3432\"],\n  [\"79\u00d778=6162\", \"49\u00d739=1911\"],\n  [\"81\u00d797=7857\", \"18\u00d776=1368\"],\n  [\"98\u00d776=7448\", \"30\u00d795=2850\"],\n  [\"33\u00d772=2376\", \"60\u00d756=3360\"],\n  [\"88\u00d762=5456\", \"27\u00d766=1782\"],\n  [\"16\u00d730=480\", \"40\u00d743=1720\"],\n  [\"39\u00d775=2925\", \"40\u00d725=1000\"],\n  [\"62\u00d727=1674\", \"84\u00d745=3780\"],\n  [\"90\u00d731=2790\", \"81\u00d765=5265\"],\n  [\"39\u00d768=2652\", \"55\u00d774=4070\"],\n  [\"41\u00d729=1189\", \"79\u00d796=7584\"],\n  [\"11\u00d756=616\", \"13\u00d750=650\"],\n  [\"28\u00d727=756\", \"71\u00d736=2556\"],\n  [\"22\u00d783=1826\", \"41\u00d715=615\"],\n  [\"62\u00d786=5332\", \"22\u00d732=704\"],\n  [\"11\u00d715=165\", \"63\u00d779=4977\"],\n  [\"70\u00d720=1400\", \"79\u00d724=1896\"],\n  [\"91\u00d715=1365\", \"42\u00d745=1890\"],\n  [\"83\u00d740=3320\", \"48\u00d782=3936\"],\n];\n\nfor (const [oldText, newText] of pairs) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length > 0) {\n    results.items[0].insertText(newText, \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Each old equation string is unique in the document, so a plain\n# Find/Replace (one occurrence each) safely retargets every cell.\n$pairs = @(\n    @(\"13\u00d721=273\", \"42\u00d729=1218\"),\n    @(\"72\u00d770=5040\", \"41\u00d714=574\"),\n    @(\"86\u00d735=3010\", \"20\u00d790=1800\"),\n    @(\"24\u00d739=936\", \"11\u00d714=154\"),\n    @(\"24\u00d795=2280\", \"47\u00d762=2914\"),\n    @(\"88\u00d757=5016\", \"78\u00d744=3432\"),\n    @(\"79\u00d778=6162\", \"49\u00d739=1911\"),\n    @(\"81\u00d797=7857\", \"18\u00d776=1368\"),\n    @(\"98\u00d776=7448\", \"30\u00d795=2850\"),\n    @(\"33\u00d772=2376\", \"60\u00d756=3360\"),\n    @(\"88\u00d762=5456\", \"27\u00d766=1782\"),\n    @(\"16\u00d730=480\", \"40\u00d743=1720\"),\n    @(\"39\u00d775=2925\", \"40\u00d725=1000\"),\n    @(\"62\u00d727=1674\", \"84\u00d745=3780\"),\n    @(\"90\u00d731=2790\", \"81\u00d765=5265\"),\n    @(\"39\u00d768=2652\", \"55\u00d774=4070\"),\n    @(\"41\u00d729=1189\", \"79\u00d796=7584\"),\n    @(\"11\u00d756=616\", \"13\u00d750=650\"),\n    @(\"28\u00d727=756\", \"71\u00d736=2556\"),\n    @(\"22\u00d783=1826\", \"41\u00d715=615\"),\n    @(\"62\u00d786=5332\", \"22\u00d732=704\"),\n    @(\"11\u00d715=165\", \"63\u00d779=4977\"),\n    @(\"70\u00d720=1400\", \"79\u00d724=1896\"),\n    @(\"91\u00d715=1365\", \"42\u00d745=1890\"),\n    @(\"83\u00d740=3320\", \"48\u00d782=3936\"),\n)\n\nforeach ($pair in $pairs) {\n    $old = $pair[0]\n    $new = $pair[1]\n    $range = $d.Content\n    $range.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n}\n"}
